# Applies the cryptos-list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.339.30"
$ws.Range("E2").Value = "  -5.88%  "
$ws.Range("D3").Value = "3.298.49"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.01"
$ws.Range("E5").Value = "  -3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.43"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.298.04"
$ws.Range("E8").Value = "  -5.11%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.36"
$ws.Range("E10").Value = "  -4.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.117"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.372"
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").Value = "3.861.27"
$ws.Range("E13").Value = "  -5.11%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "3.297.57"
$ws.Range("E16").Value = "  -5.71%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.20"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "60.530.63"
$ws.Range("E18").Value = "  -5.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.67"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.31"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("E21").Value = "  -10.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "350.18"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "3.429.23"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.18"
$ws.Range("E26").Value = "  -7.15%  "
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.152"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  -5.71%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "3.326.50"
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.58"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "157.76"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0751"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.93"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -7.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.62"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.67"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.61"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("E51").Value = "  -4.75%  "
